$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 147 -> "147（链表排序）" with a "Done" status next to it
$ws.Range("B2").Value = "147（链表排序）"
$ws.Range("C2").Value = "Done"

# 148 -> "148（归并排序）" with a "Done" status next to it
$ws.Range("B3").Value = "148（归并排序）"
$ws.Range("C3").Value = "Done"

# 75 -> "75（快排）" with a "Done" status next to it
$ws.Range("B4").Value = "75（快排）"
$ws.Range("C4").Value = "Done"

# Widen column B so the longer labels fit (matches the saved "best fit" width).
$ws.Columns.Item(2).ColumnWidth = 15.66

# Leave the cursor on C4, matching the saved selection.
[void]$ws.Range("C4").Select()
